# Add the required "experimental" boolean element to the Metadata sheet,
# and refresh the Date value to reflect the new publication timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property; its value cell (B7) needs the literal
# text "true" (not the native boolean TRUE) to match the source data,
# which was generated as plain text. Writing "true" straight into B7
# gets auto-coerced to a boolean by Excel's General-format type
# inference, so stage it in a scratch cell with a leading quote
# (forces text), copy/paste-special just the value into B7 (keeps B7's
# existing style untouched), then clean up the scratch cell.
$ws.Range("D1").Value = "'true"
$ws.Range("D1").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null
$ws.Range("D1").Clear() | Out-Null

# Row 8 = "Date" property; update its value cell (B8) to the new timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
